$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B
$ws.Cells.Item(1, 2).Value = 23
$ws.Cells.Item(2, 2).Value = 59
$ws.Cells.Item(3, 2).Value = 351
$ws.Cells.Item(5, 2).Value = 3

# Append new rows 12 and 13 with large numeric-looking IDs kept as text
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "614771330285305856"
$ws.Cells.Item(12, 1).ClearFormats()
$ws.Cells.Item(12, 2).Value = 116

$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "415445824747864064"
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(13, 2).Value = 11
